$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from Sheet1 to Raw_Data
$ws.Name = "Raw_Data"

# Fill in the additional data columns (Tested, Infected, Recovery, Deaths)
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 14

$ws.Range("C3").Value = 21
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = 11

$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 33
$ws.Range("F4").Value = 11

$ws.Range("C5").Value = 41
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 44

$ws.Range("C6").Value = 51
$ws.Range("D6").Value = 11
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = 54

# Update the selected cell to F6
$ws.Range("F6").Select()
